# Regenerate merged AHB files
# - Rename header row labels from *_old / *_new to *_FV2310 / *_FV2404
# - Add a structured Table (ListObject) over A1:U82
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Rename the header strings -------------------------------------
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J carry the "_old" -> "_FV2310" header names
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $baseNames[$i] + "_FV2310"
}

# Column K is "diff" and stays untouched.

# Columns L-U carry the "_new" -> "_FV2404" header names
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 12)
    $cell.Value = $baseNames[$i] + "_FV2404"
}

# --- 2. Turn the range into an Excel Table (ListObject) ----------------
$ws.Activate()
$ws.Range("A1").Select()
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U82"), $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row -------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
